$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): two new columns ---
#   G1 "Trường học hiện tại *"
#   H1 "Lớp học hiện tại *"
# Values are written in the same order the target workbook's shared-string
# table was built in: G1, H1, then (row 2) H2, G2.
$ws.Range("G1").Value = "Trường học hiện tại *"
$ws.Range("H1").Value = "Lớp học hiện tại *"

$ws.Range("H2").Value = "Lớp 11"
$ws.Range("G2").Value = "THPT Chuyên Nguyễn Du"

# Make the trailing "*" red & bold, like the other required-field headers
# ("Trường học hiện tại *" is 21 chars, "Lớp học hiện tại *" is 18 chars)
$g1Star = $ws.Range("G1").Characters(21, 1)
$g1Star.Font.Bold = $true
$g1Star.Font.Color = 255

$h1Star = $ws.Range("H1").Characters(18, 1)
$h1Star.Font.Bold = $true
$h1Star.Font.Color = 255

# --- Reuse existing cell styles so no new style indexes are created ---
# Header style (bold font) F1 -> G1:H1
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1:H1").PasteSpecial(-4122) | Out-Null

# Data-row text style F2 -> G2 (H2 keeps the plain/default style, like the diff)
$ws.Range("F2").Copy() | Out-Null
$ws.Range("G2").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- Column widths for the new columns (closest values this engine can store) ---
$ws.Columns.Item(7).ColumnWidth = 21.55
$ws.Columns.Item(8).ColumnWidth = 15.85
